# Applies the updated cryptos list values cell-by-cell, matching the
# source diff exactly (prices in column D, 1h volume % in column E,
# plus a coin-name/link swap between rows 40 and 41).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as text (e.g. "27.821.17", "1.008"). Excel would
# otherwise auto-detect plain decimal-looking values as numbers, which
# would change their stored representation, so force text for the whole
# affected range first, write the values, then restore the original style.
$dRange = $ws.Range("D2:D51")
$origDStyle = $dRange.Style
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.821.17"
$ws.Range("E2").Value = "  +1.71%  "

$ws.Range("D3").Value = "1.887.42"
$ws.Range("E3").Value = "  +1.57%  "

$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.55%  "

$ws.Range("D5").Value = "334.07"
$ws.Range("E5").Value = "  +1.77%  "

$ws.Range("E6").Value = "  +0.56%  "

$ws.Range("D7").Value = "0.4714"
$ws.Range("E7").Value = "  +2.26%  "

$ws.Range("D8").Value = "0.3934"
$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("D9").Value = "47.62"
$ws.Range("E9").Value = "  +2.51%  "

$ws.Range("D10").Value = "0.08084"
$ws.Range("E10").Value = "  +1.69%  "

$ws.Range("D11").Value = "1.028"
$ws.Range("E11").Value = "  +1.53%  "

$ws.Range("D12").Value = "22.23"
$ws.Range("E12").Value = "  +3.45%  "

$ws.Range("D13").Value = "1.883.38"
$ws.Range("E13").Value = "  +1.42%  "

$ws.Range("D14").Value = "5.984"
$ws.Range("E14").Value = "  +0.84%  "

$ws.Range("D15").Value = "7.135"
$ws.Range("E15").Value = "  -0.20%  "

$ws.Range("D16").Value = "1.011"
$ws.Range("E16").Value = "  +0.66%  "

$ws.Range("D17").Value = "0.06770"
$ws.Range("E17").Value = "  +2.81%  "

$ws.Range("D18").Value = "87.33"
$ws.Range("E18").Value = "  +1.23%  "

$ws.Range("D19").Value = "0.00001050"
$ws.Range("E19").Value = "  +1.82%  "

$ws.Range("D20").Value = "17.39"
$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D22").Value = "27.842.57"
$ws.Range("E22").Value = "  +1.74%  "

$ws.Range("D23").Value = "5.537"
$ws.Range("E23").Value = "  +0.95%  "

$ws.Range("D24").Value = "11.03"
$ws.Range("E24").Value = "  +1.15%  "

$ws.Range("D25").Value = "2.335"
$ws.Range("E25").Value = "  +1.45%  "

$ws.Range("D26").Value = "2.104.84"
$ws.Range("E26").Value = "  +1.58%  "

$ws.Range("D27").Value = "159.12"
$ws.Range("E27").Value = "  +3.71%  "

$ws.Range("D28").Value = "20.19"
$ws.Range("E28").Value = "  +0.37%  "

$ws.Range("D29").Value = "2.110"
$ws.Range("E29").Value = "  +2.14%  "

$ws.Range("D30").Value = "5.584"
$ws.Range("E30").Value = "  +1.84%  "

$ws.Range("D31").Value = "122.04"
$ws.Range("E31").Value = "  +0.27%  "

$ws.Range("D32").Value = "0.9834"
$ws.Range("E32").Value = "  +3.26%  "

$ws.Range("D33").Value = "0.09491"
$ws.Range("E33").Value = "  +0.67%  "

$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("D35").Value = "3.619"
$ws.Range("E35").Value = "  +0.77%  "

$ws.Range("D36").Value = "5.366"
$ws.Range("E36").Value = "  +1.92%  "

$ws.Range("D37").Value = "0.06171"
$ws.Range("E37").Value = "  +2.23%  "

$ws.Range("D38").Value = "0.02270"
$ws.Range("E38").Value = "  +1.86%  "

$ws.Range("D39").Value = "1.220"
$ws.Range("E39").Value = "  +0.73%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.6014"
$ws.Range("E40").Value = "  +1.47%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "8.065"
$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("D43").Value = "10.33"
$ws.Range("E43").Value = "  +1.42%  "

$ws.Range("D44").Value = "1.260"
$ws.Range("E44").Value = "  -1.83%  "

$ws.Range("D45").Value = "0.5730"
$ws.Range("E45").Value = "  +2.05%  "

$ws.Range("D46").Value = "12.24"
$ws.Range("E46").Value = "  +0.97%  "

$ws.Range("D47").Value = "1.949"
$ws.Range("E47").Value = "  +1.60%  "

$ws.Range("D48").Value = "3.402"
$ws.Range("E48").Value = "  +0.19%  "

$ws.Range("D49").Value = "0.06918"
$ws.Range("E49").Value = "  +2.30%  "

$ws.Range("D50").Value = "114.06"
$ws.Range("E50").Value = "  +4.43%  "

$ws.Range("D51").Value = "0.00000000303"
$ws.Range("E51").Value = "  +7.27%  "

$dRange.Style = $origDStyle

